# Refresh the cryptocurrency price list: price column D and the 1h
# volume-change percentage column E get newly scraped figures, matching
# the scheduled "Updated cryptos list ... with GitHub Actions" job.
# Rows 37/38 additionally swap: BinanceUSD now ranks above LidoDAOToken,
# so their name/link/price/change values are exchanged between the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume/1h change) hold numeric-looking text
# (e.g. "1.00", "36.448.97", "0.120") that must stay literal text rather
# than being re-interpreted as numbers (which would drop trailing zeros
# or mangle the "thousands.dot" notation). Mark the whole data range as
# Text before writing the values so they round-trip exactly.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{ Cell = 'D2'; Value = '36.448.97' },
    @{ Cell = 'E2'; Value = '  -2.88%  ' },
    @{ Cell = 'D3'; Value = '1.981.51' },
    @{ Cell = 'E3'; Value = '  -3.56%  ' },
    @{ Cell = 'E4'; Value = '  +0.26%  ' },
    @{ Cell = 'D5'; Value = '244.51' },
    @{ Cell = 'E5'; Value = '  -3.13%  ' },
    @{ Cell = 'D6'; Value = '0.627' },
    @{ Cell = 'E6'; Value = '  -3.52%  ' },
    @{ Cell = 'D7'; Value = '58.63' },
    @{ Cell = 'E7'; Value = '  -11.63%  ' },
    @{ Cell = 'E8'; Value = '  +0.09%  ' },
    @{ Cell = 'D9'; Value = '0.377' },
    @{ Cell = 'E9'; Value = '  -1.34%  ' },
    @{ Cell = 'D10'; Value = '57.54' },
    @{ Cell = 'E10'; Value = '  -4.07%  ' },
    @{ Cell = 'D11'; Value = '0.0817' },
    @{ Cell = 'E11'; Value = '  +6.35%  ' },
    @{ Cell = 'E12'; Value = '  -0.88%  ' },
    @{ Cell = 'D13'; Value = '23.78' },
    @{ Cell = 'E13'; Value = '  +11.20%  ' },
    @{ Cell = 'D14'; Value = '0.863' },
    @{ Cell = 'E14'; Value = '  -4.89%  ' },
    @{ Cell = 'D15'; Value = '14.01' },
    @{ Cell = 'E15'; Value = '  -6.29%  ' },
    @{ Cell = 'D16'; Value = '2.275.35' },
    @{ Cell = 'E16'; Value = '  -3.46%  ' },
    @{ Cell = 'D17'; Value = '5.45' },
    @{ Cell = 'E17'; Value = '  -2.39%  ' },
    @{ Cell = 'D18'; Value = '1.981.75' },
    @{ Cell = 'E18'; Value = '  -3.52%  ' },
    @{ Cell = 'D19'; Value = '36.407.81' },
    @{ Cell = 'E19'; Value = '  -2.48%  ' },
    @{ Cell = 'D20'; Value = '70.77' },
    @{ Cell = 'E20'; Value = '  -4.01%  ' },
    @{ Cell = 'D21'; Value = '0.0₃0862' },
    @{ Cell = 'E21'; Value = '  -1.70%  ' },
    @{ Cell = 'D22'; Value = '5.33' },
    @{ Cell = 'E22'; Value = '  -2.25%  ' },
    @{ Cell = 'D23'; Value = '234.83' },
    @{ Cell = 'E23'; Value = '  -2.25%  ' },
    @{ Cell = 'E24'; Value = '  +0.11%  ' },
    @{ Cell = 'D25'; Value = '2.61' },
    @{ Cell = 'E25'; Value = '  -1.48%  ' },
    @{ Cell = 'E26'; Value = '  -3.56%  ' },
    @{ Cell = 'D27'; Value = '10.15' },
    @{ Cell = 'E27'; Value = '  +3.86%  ' },
    @{ Cell = 'D28'; Value = '161.86' },
    @{ Cell = 'E28'; Value = '  +0.91%  ' },
    @{ Cell = 'D29'; Value = '19.85' },
    @{ Cell = 'E29'; Value = '  -1.02%  ' },
    @{ Cell = 'D30'; Value = '0.127' },
    @{ Cell = 'E30'; Value = '  +11.57%  ' },
    @{ Cell = 'D31'; Value = '0.120' },
    @{ Cell = 'E31'; Value = '  -1.55%  ' },
    @{ Cell = 'E32'; Value = '  +0.59%  ' },
    @{ Cell = 'D33'; Value = '4.91' },
    @{ Cell = 'E33'; Value = '  -6.70%  ' },
    @{ Cell = 'D34'; Value = '0.0632' },
    @{ Cell = 'E34'; Value = '  +1.51%  ' },
    @{ Cell = 'D35'; Value = '4.43' },
    @{ Cell = 'E35'; Value = '  -7.06%  ' },
    @{ Cell = 'D36'; Value = '6.29' },
    @{ Cell = 'E36'; Value = '  +2.47%  ' },
    @{ Cell = 'B37'; Value = 'BinanceUSD' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Cell = 'D37'; Value = '1.00' },
    @{ Cell = 'E37'; Value = '  +0.34%  ' },
    @{ Cell = 'B38'; Value = 'LidoDAOToken' },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Cell = 'D38'; Value = '2.27' },
    @{ Cell = 'E38'; Value = '  -7.26%  ' },
    @{ Cell = 'E39'; Value = '  -4.45%  ' },
    @{ Cell = 'D40'; Value = '3.10' },
    @{ Cell = 'E40'; Value = '  +2.62%  ' },
    @{ Cell = 'E41'; Value = '  +0.98%  ' },
    @{ Cell = 'D42'; Value = '0.0966' },
    @{ Cell = 'E42'; Value = '  -6.24%  ' },
    @{ Cell = 'E43'; Value = '  -3.44%  ' },
    @{ Cell = 'D44'; Value = '0.0214' },
    @{ Cell = 'E44'; Value = '  -2.53%  ' },
    @{ Cell = 'E45'; Value = '  -4.66%  ' },
    @{ Cell = 'D46'; Value = '16.31' },
    @{ Cell = 'E46'; Value = '  -3.88%  ' },
    @{ Cell = 'D47'; Value = '92.69' },
    @{ Cell = 'E47'; Value = '  -3.03%  ' },
    @{ Cell = 'D48'; Value = '7.59' },
    @{ Cell = 'E48'; Value = '  -4.77%  ' },
    @{ Cell = 'D49'; Value = '1.375.75' },
    @{ Cell = 'E49'; Value = '  -3.18%  ' },
    @{ Cell = 'D50'; Value = '2.86' },
    @{ Cell = 'E50'; Value = '  -3.26%  ' },
    @{ Cell = 'D51'; Value = '45.14' },
    @{ Cell = 'E51'; Value = '  -3.16%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
